$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 16-20) held date-formatted "null-like" sentinel values.
# Replace them with plain sequential numbers and drop the date number format
# so the cells fall back to the default/general style.
$ws.Range("C16").Value2 = 1
$ws.Range("C16").ClearFormats()

$ws.Range("C17").Value2 = 2
$ws.Range("C17").ClearFormats()

$ws.Range("C18").Value2 = 3
$ws.Range("C18").ClearFormats()

$ws.Range("C19").Value2 = 4
$ws.Range("C19").ClearFormats()

$ws.Range("C20").Value2 = 5
$ws.Range("C20").ClearFormats()

# D16/D17 contain text that looks like a formula ("=TRUE AND FALSE" / "=TRUE OR
# FALSE"). Re-enter them with a leading apostrophe so Excel stores them as a
# quote-prefixed literal string (quotePrefix style) instead of parsing them as
# formulas.
$ws.Range("D16").Value2 = "'=TRUE AND FALSE"
$ws.Range("D17").Value2 = "'=TRUE OR FALSE"

# Leave the selection where the author's last edit landed.
$ws.Range("D21").Select() | Out-Null
